# Update cryptos list: Price (D) and Volume(1h) (E) values, and two row swaps (B/C/D/E).
# Source: scheduled GitHub Actions refresh of coinranking.com data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.411.02"
$ws.Range("E2").Value = "  -3.86%  "

$ws.Range("D3").Value = "2.246.41"
$ws.Range("E3").Value = "  -5.37%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "494.82"
$ws.Range("E5").Value = "  -2.76%  "

$ws.Range("D6").Value = "127.37"
$ws.Range("E6").Value = "  -4.63%  "

$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("E8").Value = "  -2.72%  "

$ws.Range("D9").Value = "2.294.19"

$ws.Range("D10").Value = "'0.0950"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.77%  "

$ws.Range("E11").Value = "  +0.84%  "

$ws.Range("E12").Value = "  -0.46%  "

$ws.Range("E13").Value = "  -4.31%  "

$ws.Range("D14").Value = "2.695.81"
$ws.Range("E14").Value = "  -3.71%  "

$ws.Range("E15").Value = "  -1.04%  "

$ws.Range("D16").Value = "54.376.66"
$ws.Range("E16").Value = "  -3.78%  "

$ws.Range("E17").Value = "  -3.37%  "

$ws.Range("D18").Value = "2.290.21"
$ws.Range("E18").Value = "  -4.16%  "

$ws.Range("D19").Value = "10.01"
$ws.Range("E19").Value = "  -0.55%  "

$ws.Range("D20").Value = "4.06"
$ws.Range("E20").Value = "  -0.67%  "

$ws.Range("D21").Value = "'304.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.49%  "

$ws.Range("D22").Value = "6.48"
$ws.Range("E22").Value = "  +3.18%  "

$ws.Range("E23").Value = "  +0.36%  "

$ws.Range("D24").Value = "5.34"
$ws.Range("E24").Value = "  -3.67%  "

$ws.Range("D25").Value = "63.45"
$ws.Range("E25").Value = "  -3.50%  "

$ws.Range("E26").Value = "  +1.15%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.152"
$ws.Range("E27").Value = "  +2.24%  "

$ws.Range("B28").Value = "Polygon"
$ws.Range("C28").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D28").Value = "0.374"
$ws.Range("E28").Value = "  -0.90%  "

$ws.Range("D29").Value = "2.397.86"
$ws.Range("E29").Value = "  -3.51%  "

$ws.Range("D30").Value = "'7.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.69%  "

$ws.Range("D31").Value = "170.01"
$ws.Range("E31").Value = "  -1.00%  "

$ws.Range("E32").Value = "  -3.34%  "

$ws.Range("D33").Value = "0.0₃0688"
$ws.Range("E33").Value = "  -4.97%  "

$ws.Range("E34").Value = "  -0.38%  "

$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "1.08"
$ws.Range("E36").Value = "  -2.98%  "

$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "0.992"
$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("E38").Value = "  -1.04%  "

$ws.Range("E39").Value = "  -0.61%  "

$ws.Range("D40").Value = "0.869"
$ws.Range("E40").Value = "  -1.81%  "

$ws.Range("D41").Value = "3.65"
$ws.Range("E41").Value = "  -4.03%  "

$ws.Range("D42").Value = "35.53"
$ws.Range("E42").Value = "  -2.86%  "

$ws.Range("D43").Value = "0.376"
$ws.Range("E43").Value = "  -0.98%  "

$ws.Range("E44").Value = "  -2.86%  "

$ws.Range("D45").Value = "130.02"
$ws.Range("E45").Value = "  +2.19%  "

$ws.Range("E46").Value = "  -2.05%  "

$ws.Range("D47").Value = "4.94"
$ws.Range("E47").Value = "  -2.26%  "

$ws.Range("D48").Value = "0.0895"
$ws.Range("E48").Value = "  -0.92%  "

$ws.Range("D49").Value = "0.548"
$ws.Range("E49").Value = "  -2.71%  "

$ws.Range("D50").Value = "241.86"
$ws.Range("E50").Value = "  -2.70%  "

$ws.Range("D51").Value = "0.0481"
$ws.Range("E51").Value = "  -1.30%  "
